$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F3").Value = 1
$ws.Range("G3").Value = "Zjištení informací pro implementaci - frameworky"

$ws.Range("F6").Value = 1.5
$ws.Range("G6").Value = "SWOT analýza - obecné informace a vytvoření"

$ws.Range("F9").Value = 0.5
$ws.Range("G9").Value = "Přidání relevantních informací z ostatních zdrojů do SWOT"

$ws.Range("F13").Value = 0.75
$ws.Range("G13").Value = "Obecné požadavky - přepracování"

$ws.Range("F15").Value = 1.5
$ws.Range("G15").Value = "GitHub - vytvoření, nastavení, pozvánky"

$ws.Range("F16").Value = 1
$ws.Range("G16").Value = "GitHub - asistence pro ostatní"

$ws.Range("F17").Value = 6
$ws.Range("G17").Value = "Tvorba BPM"

[void]$ws.Range("G18").Select()
